$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.618.66"
$ws.Range("E2").Value = "  +3.62%  "
$ws.Range("D3").Value = "2.254.30"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.01%  "
$st = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.12"
$ws.Range("D5").Style = $st
$ws.Range("E5").Value = "  +2.65%  "
$st = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.15"
$ws.Range("D6").Style = $st
$ws.Range("E6").Value = "  +3.88%  "
$st = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.529"
$ws.Range("D7").Style = $st
$ws.Range("E7").Value = "  +3.29%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.65%  "
$st = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.07"
$ws.Range("D10").Style = $st
$ws.Range("E10").Value = "  +3.70%  "
$st = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.94"
$ws.Range("D11").Style = $st
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "2.603.32"
$ws.Range("E15").Value = "  +1.96%  "
$st = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.13"
$ws.Range("D16").Style = $st
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "2.253.19"
$ws.Range("E17").Value = "  +2.06%  "
$st = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.758"
$ws.Range("D18").Style = $st
$ws.Range("E18").Value = "  +2.85%  "
$ws.Range("D19").Value = "41.533.00"
$ws.Range("E19").Value = "  +3.66%  "
$st = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.33"
$ws.Range("D20").Style = $st
$ws.Range("E20").Value = "  +9.10%  "
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +1.59%  "
$st = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.89"
$ws.Range("D22").Style = $st
$ws.Range("E22").Value = "  +2.31%  "
$st = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.55"
$ws.Range("D23").Style = $st
$ws.Range("E23").Value = "  +1.46%  "
$st = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.07"
$ws.Range("D24").Style = $st
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E25").Value = "  +3.76%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  +4.81%  "
$st = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.90"
$ws.Range("D28").Style = $st
$st = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.49"
$ws.Range("D29").Style = $st
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("E30").Value = "  -0.73%  "
$st = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.29"
$ws.Range("D31").Style = $st
$ws.Range("E31").Value = "  +2.67%  "
$st = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.16"
$ws.Range("D32").Style = $st
$ws.Range("E32").Value = "  +6.48%  "
$ws.Range("E33").Value = "  +0.02%  "
$st = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.14"
$ws.Range("D34").Style = $st
$ws.Range("E34").Value = "  +3.62%  "
$st = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0740"
$ws.Range("D35").Style = $st
$ws.Range("E35").Value = "  +3.54%  "
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("E38").Value = "  +2.42%  "
$st = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.59"
$ws.Range("D39").Style = $st
$ws.Range("E39").Value = "  +5.50%  "
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("E41").Value = "  +2.47%  "
$st = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("D42").Style = $st
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "2.049.75"
$ws.Range("E43").Value = "  -1.24%  "
$st = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.64"
$ws.Range("D44").Style = $st
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("E45").Value = "  +2.11%  "
$st = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.15"
$ws.Range("D46").Style = $st
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("E47").Value = "  +6.13%  "
$st = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.84"
$ws.Range("D48").Style = $st
$ws.Range("E48").Value = "  +1.34%  "
$st = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("D49").Style = $st
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("E50").Value = "  +2.65%  "
$st = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.53"
$ws.Range("D51").Style = $st
$ws.Range("E51").Value = "  +6.91%  "
